$d = $word.ActiveDocument

# 1. "JAVA" -> "JAVA Desktop"
$d.Content.Find.Execute("JAVA", $false, $false, $false, $false, $false,
                         $true, 1, $false, "JAVA Desktop", 2)

# 2. "Eclipse" -> "Netbeans"
$d.Content.Find.Execute("Eclipse", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Netbeans", 2)

$d.Save()
